# Update "想去人数" (want-to-go count) figures across the sheets to reflect
# the newly generated output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 198
$ws.Range("F3").Value = 411
$ws.Range("F4").Value = 1159
$ws.Range("F5").Value = 43
$ws.Range("F6").Value = 72
$ws.Range("F7").Value = 31
$ws.Range("F8").Value = 1071
$ws.Range("F9").Value = 523
$ws.Range("F10").Value = 355
$ws.Range("F11").Value = 426
$ws.Range("F17").Value = 504
$ws.Range("F18").Value = 1455
$ws.Range("F19").Value = 5658
$ws.Range("F21").Value = 1575
$ws.Range("F22").Value = 374
$ws.Range("F23").Value = 33
$ws.Range("F25").Value = 5149
$ws.Range("F26").Value = 123
$ws.Range("F28").Value = 1512
$ws.Range("F29").Value = 17
$ws.Range("F31").Value = 657
$ws.Range("F32").Value = 78
$ws.Range("F33").Value = 68
$ws.Range("F34").Value = 3801

# --- 演出 (Performances) sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 15
$ws.Range("F5").Value = 151
$ws.Range("F8").Value = 120

# --- 本地生活 (Local life) sheet ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 582

# --- 全部类型 (All types) sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 582
$ws.Range("F5").Value = 198
$ws.Range("F6").Value = 411
$ws.Range("F7").Value = 1159
$ws.Range("F8").Value = 43
$ws.Range("F9").Value = 72
$ws.Range("F10").Value = 31
$ws.Range("F11").Value = 1071
$ws.Range("F12").Value = 355
$ws.Range("F13").Value = 426
$ws.Range("F21").Value = 504
$ws.Range("F22").Value = 1455
$ws.Range("F23").Value = 5658
$ws.Range("F25").Value = 1575
$ws.Range("F28").Value = 374
$ws.Range("F31").Value = 5149
$ws.Range("F32").Value = 123
$ws.Range("F34").Value = 1512
$ws.Range("F35").Value = 17
$ws.Range("F37").Value = 657
$ws.Range("F38").Value = 78
$ws.Range("F44").Value = 68
$ws.Range("F46").Value = 3801
